$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13 currently has only the index (A13) filled in; the rest of the
# row (B13:G13) is blank. Fill it in with the new client record, matching
# the style already used by the row above (row 12) for columns F and G.
$ws.Range("F12:G12").Copy()
$ws.Range("F13:G13").PasteSpecial(-4122)

$ws.Range("B13").Value = "ABRAÃO MOREIRA"
$ws.Range("C13").Value = "c8ea0d0ad755b73242ba6f43e8c23c6f"
$ws.Range("D13").Value = 44833.0
$ws.Range("E13").Value = 365.0
$ws.Range("F13").Value = "-"
$ws.Range("G13").Value = "VENDA 08 (29/09)"
